# Update iems-connector-test-data.xlsx test-id labels and two condition values
# to match renumbered/renamed test variants, per commit:
#   "update iems-api-engine-test-data.xlsx and iems-connector-test-data.xlsx"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("getConceptModelDataByCondition")

# --- column A (test-id) renames -------------------------------------------------
$ws.Range("A5").Value2  = "iems-connector-test-3-var2"
$ws.Range("A6").Value2  = "iems-connector-test-3-var3"
$ws.Range("A7").Value2  = "iems-connector-test-3-var4"
$ws.Range("A8").Value2  = "iems-connector-test-3-var5"

$ws.Range("A11").Value2 = "iems-connector-test-6-var2"
$ws.Range("A12").Value2 = "iems-connector-test-6-var3"

$ws.Range("A15").Value2 = "iems-connector-test-8-var2"
$ws.Range("A16").Value2 = "iems-connector-test-8-var3"
$ws.Range("A17").Value2 = "iems-connector-test-8-var4"
$ws.Range("A18").Value2 = "iems-connector-test-8-var5"

$ws.Range("A20").Value2 = "iems-connector-test-9-var2"

$ws.Range("A22").Value2 = "iems-connector-test-10-var2"
$ws.Range("A23").Value2 = "iems-connector-test-10-var3"
$ws.Range("A24").Value2 = "iems-connector-test-10-var4"
$ws.Range("A25").Value2 = "iems-connector-test-10-var5"
$ws.Range("A26").Value2 = "iems-connector-test-10-var6"
$ws.Range("A27").Value2 = "iems-connector-test-10-var7"
$ws.Range("A28").Value2 = "iems-connector-test-10-var8"
$ws.Range("A29").Value2 = "iems-connector-test-10-var9"
$ws.Range("A30").Value2 = "iems-connector-test-10-var10"

$ws.Range("A31").Value2 = "iems-connector-test-mysql-string-1"
$ws.Range("A32").Value2 = "iems-connector-test-mysql-time-1"
$ws.Range("A33").Value2 = "iems-connector-test-mysql-time-2"
$ws.Range("A34").Value2 = "iems-connector-test-mysql-time-3"
$ws.Range("A35").Value2 = "iems-connector-test-mysql-time-4"

# --- column C (condition) value updates -----------------------------------------
$ws.Range("C31").Value2 = "runid='20220227170207_2_f0d1366c-97ab-11ec-8dcf-0242ac120005'"
$ws.Range("C36").Value2 = "id='47'"

# --- restore view: frozen pane top-left cell and active selection ---------------
$ws.Range("B45").Select()
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 13
